$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, shifting existing rows 40-65 down to 41-66
$ws.Rows.Item(40).Insert()

# Populate the new row 40 with the new record's data
$ws.Cells.Item(40, 1).Value = 1
$ws.Cells.Item(40, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(40, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(40, 4).Value = 44658
$ws.Cells.Item(40, 5).Value = 15
$ws.Cells.Item(40, 6).Value = "Fruta"
$ws.Cells.Item(40, 7).Value = 100103
$ws.Cells.Item(40, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(40, 9).Value = 100103006
$ws.Cells.Item(40, 10).Value = "Nectarín"
$ws.Cells.Item(40, 11).Value = "Artic Snow"
$ws.Cells.Item(40, 12).Value = "Segunda"
$ws.Cells.Item(40, 13).Value = 250
$ws.Cells.Item(40, 14).Value = 18000
$ws.Cells.Item(40, 15).Value = 19000
$ws.Cells.Item(40, 16).Value = 18500
$ws.Cells.Item(40, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(40, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(40, 19).Value = 1028
$ws.Cells.Item(40, 20).Value = 18
